# Azure Communication Services port-in form:
# Title run "Letter of Authorization - Geographic Number Porting"
# loses the word "Geographic " and, in the process, the remaining
# text ends up split across three runs that all share the same
# (bold, 24-half-pt) run formatting:
#   "Letter of Authorization -" / " " / "Number Porting"

$d = $word.ActiveDocument

$enDash = [char]0x2013

# 1) Remove the word "Geographic " (with its trailing space) from the
#    heading so the text reads "Letter of Authorization - Number Porting".
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Text = "Geographic "
$found = $find.Find.Execute()
if ($found) {
    $find.Text = ""
}

# 2) Re-locate the (now shortened) heading text and split it into three
#    runs: "Letter of Authorization -", " " and "Number Porting" - all
#    keeping identical formatting (bold, 24 half-points).
$heading = $d.Content
$heading.Find.ClearFormatting()
$heading.Find.Text = "Letter of Authorization " + $enDash + " Number Porting"
$foundHeading = $heading.Find.Execute()

if ($foundHeading) {
    $start = $heading.Start

    # The separating space becomes its own run, splitting the text into
    # "Letter of Authorization -" / " " / "Number Porting". Nudging its
    # font size away and back forces Word to keep it (and therefore its
    # neighbours) as distinct runs instead of re-merging them, even
    # though all three end up with identical (bold, 24 half-pt) rPr.
    $partB = $d.Range($start + 25, $start + 26)      # " "
    $partB.Font.Size = 24
    $partB.Font.Size = 12
}
